# Swap the order of recorder names in column G ("Recorded By") on the
# "Session Analysis Results" sheet: cells currently reading
#   "System, dnasr281@gmail.com"
# should become
#   "dnasr281@gmail.com, System"
# Any other cell values (e.g. a single recorder name) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2
    if ($val -eq $oldValue) {
        $cell.Value = $newValue
    }
}
